$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.729797666666667
$ws.Range("H2").Value = 5.189393000000001
$ws.Range("I2").Value = 0.06436583050179444
$ws.Range("J2").Value = 0.06436583050179444
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.036942
$ws.Range("N2").Value = 0.110826
$ws.Range("O2").Value = 0.02099032928903418
$ws.Range("P2").Value = 0.02099032928903418
$ws.Range("Q2").Value = 0.063902185402
$ws.Range("R2").Value = 0.575119668618
$ws.Range("S2").Value = 0.001351059977194825
$ws.Range("T2").Value = 0.001351059977194825

# Row 3
$ws.Range("G3").Value = 1.729797666666667
$ws.Range("H3").Value = 5.189393000000001
$ws.Range("I3").Value = 0.06436583050179444
$ws.Range("J3").Value = 0.06436583050179444
$ws.Range("O3").Value = 0.5358731102718634
$ws.Range("P3").Value = 0.5358731102718634
$ws.Range("Q3").Value = 1.631392360406111
$ws.Range("R3").Value = 14.682531243655
$ws.Range("S3").Value = 0.03449191778622816
$ws.Range("T3").Value = 0.03449191778622816

# Row 4
$ws.Range("G4").Value = 1.729797666666667
$ws.Range("H4").Value = 5.189393000000001
$ws.Range("I4").Value = 0.06436583050179444
$ws.Range("J4").Value = 0.06436583050179444
$ws.Range("O4").Value = 0.4431365604391025
$ws.Range("P4").Value = 0.4431365604391026
$ws.Range("Q4").Value = 1.349068623634111
$ws.Range("R4").Value = 12.141617612707
$ws.Range("S4").Value = 0.02852285273837146
$ws.Range("T4").Value = 0.02852285273837147

# Row 5
$ws.Range("I5").Value = 0.2200595722726403
$ws.Range("J5").Value = 0.2200595722726403
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.036942
$ws.Range("N5").Value = 0.110826
$ws.Range("O5").Value = 0.02099032928903418
$ws.Range("P5").Value = 0.02099032928903418
$ws.Range("Q5").Value = 0.218474421556
$ws.Range("R5").Value = 1.966269794004
$ws.Range("S5").Value = 0.004619122885206735
$ws.Range("T5").Value = 0.004619122885206736

# Row 6
$ws.Range("I6").Value = 0.2200595722726403
$ws.Range("J6").Value = 0.2200595722726403
$ws.Range("O6").Value = 0.5358731102718634
$ws.Range("P6").Value = 0.5358731102718634
$ws.Range("S6").Value = 0.1179240074388357
$ws.Range("T6").Value = 0.1179240074388357

# Row 7
$ws.Range("I7").Value = 0.2200595722726403
$ws.Range("J7").Value = 0.2200595722726403
$ws.Range("O7").Value = 0.4431365604391025
$ws.Range("P7").Value = 0.4431365604391026
$ws.Range("S7").Value = 0.09751644194859793
$ws.Range("T7").Value = 0.09751644194859795

# Row 8
$ws.Range("H8").Value = 57.69206699999999
$ws.Range("I8").Value = 0.7155745972255653
$ws.Range("J8").Value = 0.7155745972255653
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.036942
$ws.Range("N8").Value = 0.110826
$ws.Range("O8").Value = 0.02099032928903418
$ws.Range("P8").Value = 0.02099032928903418
$ws.Range("Q8").Value = 0.7104201130379999
$ws.Range("R8").Value = 6.393781017341999
$ws.Range("S8").Value = 0.01502014642663262
$ws.Range("T8").Value = 0.01502014642663262

# Row 9
$ws.Range("H9").Value = 57.69206699999999
$ws.Range("I9").Value = 0.7155745972255653
$ws.Range("J9").Value = 0.7155745972255653
$ws.Range("O9").Value = 0.5358731102718634
$ws.Range("P9").Value = 0.5358731102718634
$ws.Range("Q9").Value = 18.13668715393833
$ws.Range("S9").Value = 0.3834571850467995
$ws.Range("T9").Value = 0.3834571850467995

# Row 10
$ws.Range("H10").Value = 57.69206699999999
$ws.Range("I10").Value = 0.7155745972255653
$ws.Range("J10").Value = 0.7155745972255653
$ws.Range("O10").Value = 0.4431365604391025
$ws.Range("P10").Value = 0.4431365604391026
$ws.Range("S10").Value = 0.3170972657521332
$ws.Range("T10").Value = 0.3170972657521332

Write-Host "Updated TPM values applied"
